$wb = $excel.ActiveWorkbook
Write-Host $wb.Worksheets.Count
foreach ($ws in $wb.Worksheets) {
    Write-Host $ws.Name
}
